# The author placed the cursor at the very start of the document (the
# "Anforderungen Roulette" heading), typed "0", and left it there - which is
# exactly where Word's automatic "_GoBack" bookmark (last-edit-location
# marker) ends up after that edit. Reproduce both effects:
#   1. Insert a new run containing "0" right before the existing
#      "Anforderungen " run in the first paragraph (the heading).
#   2. Move the "_GoBack" bookmark so it wraps the position right after
#      that new "0" text (Word only ever keeps a single "_GoBack"
#      bookmark, so adding it here automatically removes it from its old
#      location at the end of the "Dutzend: ..." list paragraph).

$d = $word.ActiveDocument

# 1. Insert "0" at the very beginning of the document (start of heading).
$headingStart = $d.Paragraphs(1).Range.Duplicate
$headingStart.Collapse(1)
$headingStart.InsertBefore("0")

# 2. Re-create "_GoBack" right after the inserted "0" (collapsed range at
#    document position 1, i.e. right after the single "0" character).
$goBackRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
